# Clase 14 Junio 2017
# Normalize the "Sexo" column values: "Hembra" -> "Mujer", "Macho" -> "Hombre"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F13").Value = "Mujer"
$ws.Range("F16").Value = "Hombre"
